$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert the first new record as row 300 (right after the existing row 299),
# which pushes the former rows 300-329 down to 301-330.
$ws.Rows.Item(300).Insert()
$ws.Cells.Item(300,1).Value = 11
$ws.Cells.Item(300,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(300,3).Value = "Bíobío"
$ws.Cells.Item(300,4).Value = 44679
$ws.Cells.Item(300,5).Value = 8
$ws.Cells.Item(300,6).Value = 100112006
$ws.Cells.Item(300,7).Value = "Repollo"
$ws.Cells.Item(300,8).Value = "Morada(o)"
$ws.Cells.Item(300,9).Value = "Primera"
$ws.Cells.Item(300,10).Value = 1000
$ws.Cells.Item(300,11).Value = 1600
$ws.Cells.Item(300,12).Value = 1700
$ws.Cells.Item(300,13).Value = 1650
$ws.Cells.Item(300,14).Value = "`$/unidad"
$ws.Cells.Item(300,15).Value = "Región Metropolitana"
$ws.Cells.Item(300,16).Value = 1650
$ws.Cells.Item(300,17).Value = 1
$ws.Cells.Item(300,18).Value = "Hortaliza"

# Insert the second new record as row 319 (right after the record that was
# originally row 318, now sitting at row 319 because of the insert above),
# which pushes the remaining rows further down by one more.
$ws.Rows.Item(319).Insert()
$ws.Cells.Item(319,1).Value = 11
$ws.Cells.Item(319,2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(319,3).Value = "Bíobío"
$ws.Cells.Item(319,4).Value = 44629
$ws.Cells.Item(319,5).Value = 8
$ws.Cells.Item(319,6).Value = 100112006
$ws.Cells.Item(319,7).Value = "Repollo"
$ws.Cells.Item(319,8).Value = "Morada(o)"
$ws.Cells.Item(319,9).Value = "Primera"
$ws.Cells.Item(319,10).Value = 500
$ws.Cells.Item(319,11).Value = 1700
$ws.Cells.Item(319,12).Value = 1800
$ws.Cells.Item(319,13).Value = 1740
$ws.Cells.Item(319,14).Value = "`$/unidad"
$ws.Cells.Item(319,15).Value = "Región Metropolitana"
$ws.Cells.Item(319,16).Value = 1740
$ws.Cells.Item(319,17).Value = 1
$ws.Cells.Item(319,18).Value = "Hortaliza"
